# ELEC374 MP1 values - add new "Time vs Matrix Width" table broken out per
# block/tile width (TW = 2, 5, 10, 25) in columns T:X, mirroring the layout
# already used for the J:O block-width table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3): reuse the "Matrix Width" label, then one column
# per tile width. Strings are written in the same order Excel's original
# author workbook introduced them, so the shared-string table lines up.
$ws.Range("V3").Value = "TW = 5"
$ws.Range("W3").Value = "TW = 10"
$ws.Range("X3").Value = "TW = 25"
$ws.Range("U3").Value = "TW = 2"
$ws.Range("T3").Value = "Matrix Width"

# --- Data rows 4-28: Matrix Width, then timing values for TW=2,5,10,25.
$data = @(
    @(4, 100, 7.7952, 1.304096, 0.800768, 0.712032),
    @(5, 100, 6.73936, 1.280704, 0.782336, 0.698688),
    @(6, 100, 8.2264, 1.28, 0.784384, 0.69632),
    @(7, 100, 6.733728, 1.294304, 0.782336, 0.696064),
    @(8, 100, 6.734976, 1.279744, 0.78224, 1.776544),
    @(9, 250, 112.655075, 20.342943, 7.703488, 8.32224),
    @(10, 250, 83.240448, 20.771551, 7.668544, 8.295552),
    @(11, 250, 84.441277, 16.817345, 7.673856, 11.24368),
    @(12, 250, 51.855553, 20.404257, 9.7968, 8.3184),
    @(13, 250, 43.494946, 16.808672, 7.6816, 8.296768),
    @(14, 500, 498.086853, 135.974777, 67.804031, 73.433792),
    @(15, 500, 337.866943, 115.226524, 59.093346, 73.253922),
    @(16, 500, 332.337708, 73.752609, 62.290241, 45.690113),
    @(17, 500, 337.924652, 55.515327, 60.696831, 32.927776),
    @(18, 500, 334.71225, 56.073566, 35.610592, 32.799774),
    @(19, 1000, 2946.864258, 593.702515, 354.925842, 367.698761),
    @(20, 1000, 2779.178223, 439.831055, 213.968414, 224.805984),
    @(21, 1000, 2794.071289, 441.609955, 213.330078, 225.45488),
    @(22, 1000, 2798.993164, 438.207886, 210.610107, 221.807938),
    @(23, 1000, 2800.180664, 441.676147, 217.816452, 226.771652),
    @(24, 1500, 9755.364258, 1660.970825, 871.172668, 921.502136),
    @(25, 1500, 9528.395508, 1508.836426, 721.274658, 760.004944),
    @(26, 1500, 9578.418945, 1504.605713, 715.679321, 754.128174),
    @(27, 1500, 9577.771484, 1506.230713, 710.140503, 750.335266),
    @(28, 1500, 9637.646484, 1505.693237, 715.723816, 751.165894)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("T$r").Value = $entry[1]
    $ws.Range("U$r").Value = $entry[2]
    $ws.Range("V$r").Value = $entry[3]
    $ws.Range("W$r").Value = $entry[4]
    $ws.Range("X$r").Value = $entry[5]
}

# --- Cosmetic view state matching the author's last-saved window (zoom
# level and active selection).
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("M35").Select() | Out-Null

Write-Host "Edit applied."
